# feat: add 2022-Q3 data
#
# Before:  Sheet1 "总计" (totals) + Sheet2 "2022-Q2" (fund holdings for Q2)
# After:   Sheet1 "总计" (totals, +1 row for Q3) + Sheet2 "2022-Q3" (new fund
#          holdings data) + Sheet3 "2022-Q2" (the original Q2 fund holdings,
#          now moved one tab to the right, unchanged)

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item(1)     # "总计"
$wsQ2    = $wb.Worksheets.Item(2)     # "2022-Q2" (original)

# ---------------------------------------------------------------------
# 1) Duplicate the existing "2022-Q2" sheet (values + formatting) into a
#    brand-new sheet placed right after it. Doing the copy BEFORE any
#    data overwrite keeps the original Q2 data intact for the new tab.
# ---------------------------------------------------------------------
$wsQ2Copy = $wb.Worksheets.Add($null, $wsQ2)

$wsQ2.Range("B1:H1").Copy($wsQ2Copy.Range("B1:H1"))
$wsQ2.Range("A2:H2").Copy($wsQ2Copy.Range("A2:H2"))

# ---------------------------------------------------------------------
# 2) Turn the original sheet into "2022-Q3" and replace its contents with
#    the new quarter's fund-holding data. Rename the original sheet first
#    so the temporary name never collides with the copy.
# ---------------------------------------------------------------------
$wsQ2.Name = "2022-Q3"
$wsQ3 = $wsQ2

$wsQ2Copy.Name = "2022-Q2"

# Clear the old single data row so nothing lingers past row 4.
$wsQ3.Range("A2:H2").ClearContents()

# Headers (style copied from the "总计" sheet's header cell).
$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 2 + $i   # headers start at column B
    $cell = $wsQ3.Cells.Item(1, $col)
    $wsTotal.Range("B1").Copy($cell)
    $cell.Value = $headers[$i]
}

# Data rows.
$rows = @(
    @(0, "012060", "富国全球消费精选混合（QDII）A",      "2.68",  "66.08", "3.17", "0.0850", 6),
    @(1, "012061", "富国全球消费精选混合（QDII）美元现汇", "2.68",  "66.08", "3.17", "0.0850", 6),
    @(2, "012062", "富国全球消费精选混合（QDII）C",      "-2.54", "66.08", "3.17", "-0.0805", 6)
)

$r = 2
foreach ($row in $rows) {
    $cellA = $wsQ3.Cells.Item($r, 1)
    $wsTotal.Range("A2").Copy($cellA)
    $cellA.Value = $row[0]

    $cellB = $wsQ3.Cells.Item($r, 2)
    $cellB.NumberFormat = "@"
    $cellB.Value = $row[1]

    $cellC = $wsQ3.Cells.Item($r, 3)
    $cellC.Value = $row[2]

    $cellD = $wsQ3.Cells.Item($r, 4)
    $cellD.NumberFormat = "@"
    $cellD.Value = $row[3]

    $cellE = $wsQ3.Cells.Item($r, 5)
    $cellE.NumberFormat = "@"
    $cellE.Value = $row[4]

    $cellF = $wsQ3.Cells.Item($r, 6)
    $cellF.NumberFormat = "@"
    $cellF.Value = $row[5]

    $cellG = $wsQ3.Cells.Item($r, 7)
    $cellG.NumberFormat = "@"
    $cellG.Value = $row[6]

    $cellH = $wsQ3.Cells.Item($r, 8)
    $cellH.Value = $row[7]

    $r = $r + 1
}

# ---------------------------------------------------------------------
# 3) Update the "总计" sheet: push the existing Q2 summary row down to
#    row 3 and insert the new Q3 summary row above it at row 2.
# ---------------------------------------------------------------------
$wsTotal.Range("A2").Copy($wsTotal.Range("A3"))
$wsTotal.Cells.Item(3, 1).Value = 1
$wsTotal.Cells.Item(3, 2).Value = "2022-Q2"
$wsTotal.Cells.Item(3, 3).Value = 1
$wsTotal.Cells.Item(3, 4).Value = 0.19

$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q3"
$wsTotal.Cells.Item(2, 3).Value = 3
$wsTotal.Cells.Item(2, 4).Value = 0.09
